# "New Internship Page Updated"
# Rebuilds the "Students Data" sheet header row with the full internship
# column set, rewrites row 2 (Riya Ingale) with the updated internship
# details, and removes the second student row (Samiksha) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend headers from C1 out to T1 ------------------------------
$headers = @(
    "Email",
    "Department",
    "Division",
    "Year",
    "Company Name",
    "Position",
    "Domain",
    "Source",
    "skills_required",
    "Company Representative Name",
    "Company Representative Contact",
    "Start Date",
    "End Date",
    "Feedback",
    "Work Environment Rating",
    "Satisfaction",
    "Would student recommend?"
)
$col = 4
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# --- Row 2: update the internship details for Riya Ingale -----------------
$ws.Range("I2").Value = "Python Intern"
$ws.Range("J2").Value = "Web Development"
$ws.Range("K2").Value = "Self"

# O2/P2 held numeric dates styled with a custom date format; they become
# plain text strings, so drop the old style before writing the new values.
$ws.Range("O2").Style = "Normal"
$ws.Range("O2").Value = "22 /02 /21"
$ws.Range("P2").Style = "Normal"
$ws.Range("P2").Value = "21/ 03/ 21"

$ws.Range("Q2").Value = "Nice Experience"
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = "Yes"
$ws.Range("T2").Value = "Maybe"

# --- Row 3: remove the second student (Samiksha) entirely ------------------
$ws.Rows("3:3").Delete()
